$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F2").Value = 170
$ws1.Range("F3").Value = 97
$ws1.Range("F6").Value = 4976
$ws1.Range("F7").Value = 391
$ws1.Range("F8").Value = 568
$ws1.Range("F9").Value = 867
$ws1.Range("F10").Value = 803
$ws1.Range("F12").Value = 16
$ws1.Range("F13").Value = 540
$ws1.Range("F16").Value = 1626
$ws1.Range("F17").Value = 1426
$ws1.Range("F18").Value = 711
$ws1.Range("F20").Value = 175
$ws1.Range("F21").Value = 261
$ws1.Range("F22").Value = 469
$ws1.Range("F23").Value = 111
$ws1.Range("F24").Value = 1037
$ws1.Range("F27").Value = 1823
$ws1.Range("F28").Value = 150
$ws1.Range("F29").Value = 80
$ws1.Range("F30").Value = 13
$ws1.Range("F31").Value = 203
$ws1.Range("F36").Value = 255
$ws1.Range("F37").Value = 555
$ws1.Range("F39").Value = 24
$ws1.Range("F40").Value = 22

# Sheet 2
$ws2.Range("F4").Value = 135
$ws2.Range("F7").Value = 7

# Sheet 4
$ws4.Range("F3").Value = 170
$ws4.Range("F4").Value = 97
$ws4.Range("F8").Value = 4976
$ws4.Range("F9").Value = 391
$ws4.Range("F10").Value = 568
$ws4.Range("F12").Value = 135
$ws4.Range("F13").Value = 867
$ws4.Range("F14").Value = 803
$ws4.Range("F18").Value = 16
$ws4.Range("F19").Value = 540
$ws4.Range("F21").Value = 7
$ws4.Range("F23").Value = 1626
$ws4.Range("F24").Value = 1426
$ws4.Range("F25").Value = 711
$ws4.Range("F27").Value = 175
$ws4.Range("F28").Value = 261
$ws4.Range("F30").Value = 469
$ws4.Range("F31").Value = 111
$ws4.Range("F32").Value = 1037
$ws4.Range("F34").Value = 1823
$ws4.Range("F35").Value = 150
$ws4.Range("F36").Value = 80
$ws4.Range("F37").Value = 13
$ws4.Range("F38").Value = 203
$ws4.Range("F42").Value = 255
$ws4.Range("F43").Value = 555
$ws4.Range("F45").Value = 24
$ws4.Range("F46").Value = 22
